$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the very top; everything currently on the sheet
# (rows 1-28) shifts down to rows 2-29.
$ws.Rows("1:1").Insert()

# New header row describing the three columns.
$ws.Range("A1").Value = "Дата"
$ws.Range("B1").Value = "Описание"
$ws.Range("C1").Value = "Подпись"

# Header is bold; keep the date column left aligned like the rest of column A.
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4131

# New diary entry for 2024-06-25 (row shifted from 6 to 7 after the insert).
$ws.Range("B7").Value = "Изучение методов подбора сечения для балки по сортаменту"

# Entry for 2024-07-06 (shifted from row 17 to row 18) is cleared out.
$ws.Range("B18").Clear()

# Column B grew wider to fit the new longer entry; column C became the
# narrow signature column.
$ws.Columns("B:C").AutoFit()

# Restore view state.
$ws.Application.ActiveWindow.Zoom = 197
$ws.Range("B9").Select()
